$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.438.72"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "1.830.34"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.85%  "
$ws.Range("D5").Value = "'330.85"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").Value = "'0.4593"
$ws.Range("E7").Value = "  -2.55%  "
$ws.Range("D8").Value = "'0.3830"
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("D9").Value = "'46.57"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("D10").Value = "'0.07901"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").Value = "'0.9705"
$ws.Range("E11").Value = "  -3.29%  "
$ws.Range("D12").Value = "'21.08"
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("D13").Value = "1.842.86"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").Value = "'5.877"
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("D15").Value = "'7.050"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").Value = "'87.94"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "'0.06636"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").Value = "'17.16"
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("D22").Value = "27.440.14"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("D23").Value = "'5.336"
$ws.Range("E23").Value = "  -2.15%  "
$ws.Range("D24").Value = "'10.80"
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("D25").Value = "'2.300"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("D26").Value = "2.056.25"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("D27").Value = "'156.92"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "'19.41"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("D29").Value = "'2.059"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("D30").Value = "'5.265"
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("D31").Value = "'118.69"
$ws.Range("E31").Value = "  -1.92%  "
$ws.Range("D32").Value = "'0.9537"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").Value = "'0.09296"
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("D34").Value = "'3.585"
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("D35").Value = "'5.232"
$ws.Range("E35").Value = "  -1.06%  "
$ws.Range("D36").Value = "'1.312"
$ws.Range("E36").Value = "  -2.26%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02199"
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.05930"
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("D39").Value = "'8.032"
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("E40").Value = "  -4.30%  "
$ws.Range("D41").Value = "'0.5790"
$ws.Range("E41").Value = "  -1.87%  "
$ws.Range("D42").Value = "'0.1838"
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("D44").Value = "'1.274"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("D45").Value = "'0.5479"
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("D46").Value = "'11.95"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "'1.869"
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("D49").Value = "'110.22"
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("D50").Value = "'1.039"
$ws.Range("E50").Value = "  -2.15%  "
$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = "  -0.89%  "
